# Adds two new "tables" (Modulo, Exponencial) below the existing "Branchs"
# table, mirroring its layout (title row, header row, blank data row), plus
# a spacer row with wrap-text formatting between the two new tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Table 2: "Modulo" -- rows 8 (blank spacer), 9 (title), 10 (headers),
# 11 (blank values), 12 (blank wrap-text spacer row A:F)
# ---------------------------------------------------------------------

# Row 9: merged title cell, styled like the "Branchs" title (row 3).
# Merge BEFORE copying the format over so the engine doesn't redistribute
# the border across the merged cells (it only does that when merging an
# already-bordered range).
$ws.Range("B9:E9").Merge()
$ws.Range("B3:E3").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Modulo"

# Row 10: header cells, styled like row 4.
$ws.Range("B4").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("B10").Value = "cond"
$ws.Range("C10").Value = "op"
$ws.Range("D10").Value = "funct"
$ws.Range("E10").Value = "Label   "

# Row 11: blank value cells, styled like row 5 (no content).
$ws.Range("B5").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)

# Row 12: blank spacer row spanning A:F with wrap-text formatting.
$ws.Range("A12:F12").WrapText = $true

# ---------------------------------------------------------------------
# Table 3: "Exponencial" -- rows 13 (title), 14 (headers), 15 (blank
# values), 16 (blank trailing spacer)
# ---------------------------------------------------------------------

$ws.Range("B13:E13").Merge()
$ws.Range("B3:E3").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)
$ws.Range("B13").Value = "Exponencial"

$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("B14").Value = "cond"
$ws.Range("C14").Value = "op"
$ws.Range("D14").Value = "funct"
$ws.Range("E14").Value = "Label   "

$ws.Range("B5").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)

# Row 16 stays blank (trailing spacer row, default formatting) -- nothing
# to write there.

# ---------------------------------------------------------------------
# Sheet-level bookkeeping to match the edited state.
# ---------------------------------------------------------------------
$null = $ws.Range("E6").Select()
